$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.337.86"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.870.89"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'243.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.4691"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").Value = "'0.2874"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "'0.06450"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'22.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.32%  "
$ws.Range("D11").Value = "'0.07774"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "1.874.00"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'95.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "'0.7211"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "'5.126"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'278.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "30.332.24"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'12.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007524"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Value = "2.116.95"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").Value = "'5.241"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "'6.239"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").Value = "'163.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "'9.040"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").Value = "'18.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'1.874"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("D29").Value = "'1.321"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").Value = "'0.09580"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.16%  "
$ws.Range("D31").Value = "'1.468"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").Value = "'4.203"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.32%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.6885"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'2.708"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("D38").Value = "'0.01878"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.46%  "
$ws.Range("D39").Value = "'2.810"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("D40").Value = "'6.215"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'74.16"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.99%  "
$ws.Range("D42").Value = "'1.941"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").Value = "'0.4221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.93%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'0.8246"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'100.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("D47").Value = "'9.574"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "'6.920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.65%  "
$ws.Range("D50").Value = "'896.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  +0.82%  "
